$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on cells we touch so numeric-looking strings
# (e.g. "608.20", "0.0000281") are stored as text, matching the source data
# (inline strings), not auto-converted to numbers by Excel.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '70.407.52'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.66%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.530.52'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.60%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '608.20'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +4.51%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '173.31'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.47%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.618'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.88%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.520.79'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.58%  '

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.07%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.201'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +5.98%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.74'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.11%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.585'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -2.01%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '47.53'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.01%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000281'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +1.77%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.088.19'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.42%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '627.45'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -7.00%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '8.44'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -3.19%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '70.242.41'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.49%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.520.02'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.40%  '

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -2.02%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.39'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.32%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.888'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.77%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.99'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -10.72%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '15.89'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.66%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '96.62'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -1.27%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.86'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.32%  '

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.11%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.61'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.77%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.25'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.15%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '33.32'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.20%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.47'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.93%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.11'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -2.92%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.34'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.98%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.03'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.03%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '561.34'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -5.93%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '10.80'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.87%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.56'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.27%  '

$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.102'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -2.28%  '

$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '57.26'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.12%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.26%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.143'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +4.90%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0451'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.70%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.329'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -2.25%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.343.46'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.04%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.02'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +4.09%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0₃0715'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.80%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '33.17'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.76%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.65'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.55%  '

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -2.97%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '134.55'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +2.10%  '

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.64%  '
